$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9: update title and link
$ws.Range("D9").Value = "2022학년도 봄학기 지원자 지원동기"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/spring-2022-applicants-motivations/#utm_source=rss&utm_medium=rss&utm_campaign=spring-2022-applicants-motivations"

# Row 26: update title
$ws.Range("D26").Value = "2021 인공지능 경진대회 참가기"

# Row 37: update title and link
$ws.Range("D37").Value = "[Paper Review] Speech to Speech Translation"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1883&mod=document&pageid=1"
